$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new columns E-H inserted (shift old E-H content to I-L) ---
# Read the old values in E1:H1 before we overwrite them
$oldE1 = $ws.Range("E1").Value()
$oldF1 = $ws.Range("F1").Value()
$oldG1 = $ws.Range("G1").Value()
$oldH1 = $ws.Range("H1").Value()

$ws.Range("I1").Value = $oldE1
$ws.Range("J1").Value = $oldF1
$ws.Range("K1").Value = $oldG1
$ws.Range("L1").Value = $oldH1

$ws.Range("E1").Value = "Data4"
$ws.Range("F1").Value = "Data5"
$ws.Range("G1").Value = "Data6"
$ws.Range("H1").Value = "Data7"

# --- Row 2: move old E2:H2 values to I2:L2, set new D2, clear old E2:H2 ---
$oldE2 = $ws.Range("E2").Value()
$oldF2 = $ws.Range("F2").Value()
$oldG2 = $ws.Range("G2").Value()
$oldH2 = $ws.Range("H2").Value()

$ws.Range("I2").Value = $oldE2
$ws.Range("J2").Value = $oldF2
$ws.Range("K2").Value = $oldG2
$ws.Range("L2").Value = $oldH2

$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("H2").ClearContents()

$ws.Range("D2").Value = "Ricolino"

# --- New row 10 ---
$ws.Range("A10").Value = "SP_TC_52"
$ws.Range("B10").Value = "Print Preview"
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = "JINCY_SKU_0"
$ws.Range("E10").Value = "ORD REQ DATE:"
$ws.Range("F10").Value = "TOTAL       1        `$9.00"
$ws.Range("G10").Value = "DELIVERY DATE:"
$ws.Range("H10").Value = "ROUTE CODE: 4001"

# --- Remove row 13 content (B13 had style s=4, no value) ---
$ws.Range("B13").ClearContents()

# --- Column widths ---
$ws.Range("B:B").ColumnWidth = 22.28515625
$ws.Range("D:D").ColumnWidth = 15.140625
$ws.Range("E:E").ColumnWidth = 14.5703125
$ws.Range("F:F").ColumnWidth = 23.42578125
$ws.Range("G:G").ColumnWidth = 23.42578125
$ws.Range("H:H").ColumnWidth = 23.42578125

# --- Selection ---
$ws.Range("A11").Select()
